$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 651, pushing the former rows 651-760 down to 652-761.
$ws.Rows("651:651").Insert()

# Populate the newly inserted row 651 with the new weekly record.
$ws.Range("A651").Value2 = 10
$ws.Range("B651").Value2 = "Vega Modelo de Temuco"
$ws.Range("C651").Value2 = "La Araucanía"
$ws.Range("D651").Value2 = 45218
$ws.Range("E651").Value2 = 9
$ws.Range("F651").Value2 = 100112028
$ws.Range("G651").Value2 = "Sandia"
$ws.Range("H651").Value2 = "Sin especificar"
$ws.Range("I651").Value2 = "Primera"
$ws.Range("J651").Value2 = 1600
$ws.Range("K651").Value2 = 850
$ws.Range("L651").Value2 = 900
$ws.Range("M651").Value2 = 881
$ws.Range("N651").Value2 = "`$/kilo (volumen en unidades)"
$ws.Range("O651").Value2 = "Perú"
$ws.Range("P651").Value2 = 881
$ws.Range("Q651").Value2 = 1
$ws.Range("R651").Value2 = "Hortaliza"
